$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 72, shifting existing rows
# 72-87 down to 74-89 (dimension grows from A1:R87 to A1:R89).
$ws.Rows("72:73").Insert()

# New row 72: Rabanito entry with Fecha 44855 (2022-10-21) and updated
# volume/price figures.
$ws.Cells.Item(72,1).Value = 10
$ws.Cells.Item(72,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(72,3).Value = "La Araucanía"
$ws.Cells.Item(72,4).Value = 44855
$ws.Cells.Item(72,5).Value = 9
$ws.Cells.Item(72,6).Value = 300000001
$ws.Cells.Item(72,7).Value = "Rabanito"
$ws.Cells.Item(72,8).Value = "Sin especificar"
$ws.Cells.Item(72,9).Value = "Primera"
$ws.Cells.Item(72,10).Value = 40
$ws.Cells.Item(72,11).Value = 7000
$ws.Cells.Item(72,12).Value = 8000
$ws.Cells.Item(72,13).Value = 7500
$ws.Cells.Item(72,14).Value = "`$/docena de paquetes"
$ws.Cells.Item(72,15).Value = "Provincia de Cautín"
$ws.Cells.Item(72,16).Value = 625
$ws.Cells.Item(72,17).Value = 12
$ws.Cells.Item(72,18).Value = "Hortaliza"

# New row 73: Rabanito entry with Fecha 44855 (2022-10-21) sourced from
# Región Metropolitana.
$ws.Cells.Item(73,1).Value = 10
$ws.Cells.Item(73,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(73,3).Value = "La Araucanía"
$ws.Cells.Item(73,4).Value = 44855
$ws.Cells.Item(73,5).Value = 9
$ws.Cells.Item(73,6).Value = 300000001
$ws.Cells.Item(73,7).Value = "Rabanito"
$ws.Cells.Item(73,8).Value = "Sin especificar"
$ws.Cells.Item(73,9).Value = "Primera"
$ws.Cells.Item(73,10).Value = 30
$ws.Cells.Item(73,11).Value = 6000
$ws.Cells.Item(73,12).Value = 6000
$ws.Cells.Item(73,13).Value = 6000
$ws.Cells.Item(73,14).Value = "`$/docena de paquetes"
$ws.Cells.Item(73,15).Value = "Región Metropolitana"
$ws.Cells.Item(73,16).Value = 500
$ws.Cells.Item(73,17).Value = 12
$ws.Cells.Item(73,18).Value = "Hortaliza"
